$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume(1h) figures (and the two swapped rows
# for Filecoin/EnergySwap) per the latest data pull.
#
# Columns D/E hold plain text (e.g. "1.00", "4.05", "  -3.77%  ") rather
# than numeric values. Assigning a bare numeric-looking string via .Value
# would make Excel auto-coerce it to a Double (dropping the thousands-dot
# grouping / significant trailing zeros), so every such literal is given a
# leading apostrophe to force text entry, exactly like typing it by hand.

$ws.Range("D2").Value = "'63.356.71"
$ws.Range("E2").Value = "'  -3.77%  "
$ws.Range("D3").Value = "'3.123.70"
$ws.Range("E3").Value = "'  -4.82%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'559.99"
$ws.Range("E5").Value = "'  -4.48%  "
$ws.Range("D6").Value = "'161.37"
$ws.Range("E6").Value = "'  -9.48%  "
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'0.581"
$ws.Range("E8").Value = "'  -9.08%  "
$ws.Range("D9").Value = "'3.117.80"
$ws.Range("E9").Value = "'  -5.06%  "
$ws.Range("E10").Value = "'  -2.14%  "
$ws.Range("E11").Value = "'  -7.86%  "
$ws.Range("D12").Value = "'0.378"
$ws.Range("E12").Value = "'  -5.50%  "
$ws.Range("D13").Value = "'3.665.93"
$ws.Range("E13").Value = "'  -4.97%  "
$ws.Range("E14").Value = "'  -1.20%  "
$ws.Range("D15").Value = "'63.338.96"
$ws.Range("E15").Value = "'  -3.98%  "
$ws.Range("D16").Value = "'24.71"
$ws.Range("E16").Value = "'  -6.26%  "
$ws.Range("D17").Value = "'3.122.82"
$ws.Range("E17").Value = "'  -4.71%  "
$ws.Range("E18").Value = "'  -6.38%  "
$ws.Range("D19").Value = "'398.75"
$ws.Range("E19").Value = "'  -5.26%  "
$ws.Range("E20").Value = "'  -5.01%  "
$ws.Range("D21").Value = "'12.42"
$ws.Range("E21").Value = "'  -4.64%  "
$ws.Range("E22").Value = "'  -3.09%  "
$ws.Range("E23").Value = "'  +0.06%  "
$ws.Range("D24").Value = "'67.30"
$ws.Range("E24").Value = "'  -5.32%  "
$ws.Range("D25").Value = "'0.200"
$ws.Range("E25").Value = "'  -3.29%  "
$ws.Range("D26").Value = "'0.476"
$ws.Range("E26").Value = "'  -5.83%  "
$ws.Range("E27").Value = "'  -11.50%  "
$ws.Range("E28").Value = "'  -7.51%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  -0.20%  "
$ws.Range("E31").Value = "'  -6.97%  "
$ws.Range("D32").Value = "'20.85"
$ws.Range("E32").Value = "'  -5.95%  "
$ws.Range("E33").Value = "'  -5.18%  "
$ws.Range("D34").Value = "'4.76"
$ws.Range("E34").Value = "'  -7.16%  "
$ws.Range("D35").Value = "'1.10"
$ws.Range("E35").Value = "'  -7.09%  "
$ws.Range("D36").Value = "'152.42"
$ws.Range("E36").Value = "'  -3.57%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "'  -8.10%  "
$ws.Range("D38").Value = "'2.746.80"
$ws.Range("E38").Value = "'  -3.64%  "
$ws.Range("E39").Value = "'  -8.02%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.05"
$ws.Range("E40").Value = "'  -6.52%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'23.24"
$ws.Range("E41").Value = "'  -11.14%  "
$ws.Range("D42").Value = "'38.23"
$ws.Range("E42").Value = "'  -3.39%  "
$ws.Range("E43").Value = "'  -7.38%  "
$ws.Range("D44").Value = "'0.0611"
$ws.Range("E44").Value = "'  -3.98%  "
$ws.Range("D45").Value = "'5.38"
$ws.Range("E45").Value = "'  -8.32%  "
$ws.Range("E46").Value = "'  -4.40%  "
$ws.Range("D47").Value = "'20.71"
$ws.Range("E47").Value = "'  -9.19%  "
$ws.Range("E48").Value = "'  -0.09%  "
$ws.Range("D49").Value = "'279.48"
$ws.Range("E49").Value = "'  -10.09%  "
$ws.Range("E50").Value = "'  -4.65%  "
$ws.Range("E51").Value = "'  +1.07%  "
